{"js": "// Replace each arithmetic expression in the table with its updated value.\n// Every expression (e.g. \"86-67=19\") occurs exactly once in the document,\n// so a direct search-and-replace on the unique old text is unambiguous.\nconst replacements = [\n    [\"86-67=19\", \"57-16=41\"],\n    [\"52+16=68\", \"41+28=69\"],\n    [\"28+51=79\", \"49-37=12\"],\n    [\"70+22=92\", \"60-24=36\"],\n    [\"19-5=14\", \"24+70=94\"],\n    [\"88-20=68\", \"83-30=53\"],\n    [\"13-6=7\", \"16+21=37\"],\n    [\"62-5=57\", \"87-43=44\"],\n    [\"36+7=43\", \"15+16=31\"],\n    [\"86-16=70\", \"66-49=17\"],\n    [\"51-3=48\", \"16+82=98\"],\n    [\"74-52=22\", \"4+38=42\"],\n    [\"41+55=96\", \"68+24=92\"],\n    [\"54-45=9\", \"34+29=63\"],\n    [\"12-4=8\", \"52-16=36\"],\n    [\"45+35=80\", \"33+66=99\"],\n    [\"8+15=23\", \"96-65=31\"],\n    [\"72-62=10\", \"4+52=56\"],\n    [\"78-4=74\", \"98-91=7\"],\n    [\"93-65=28\", \"41-39=2\"],\n    [\"14-11=3\", \"13+68=81\"],\n    [\"27+66=93\", \"71-54=17\"],\n    [\"28+26=54\", \"89-14=75\"],\n    [\"61+6=67\", \"18+67=85\"],\n    [\"20+57=77\", \"99-59=40\"],\n    [\"18+43=61\", \"43+53=96\"],\n    [\"40-29=11\", \"13+20=33\"],\n    [\"23-17=6\", \"79+4=83\"],\n    [\"64-10=54\", \"34-34=0\"],\n    [\"78+14=92\", \"37+28=65\"],\n    [\"67+28=95\", \"20-9=11\"],\n    [\"32-5=27\", \"48+17=65\"],\n    [\"2+13=15\", \"33-2=31\"],\n    [\"86-62=24\", \"6+16=22\"],\n    [\"86-13=73\", \"23+1=24\"],\n    [\"76-36=40\", \"85-10=75\"],\n    [\"93-22=71\", \"8+32=40\"],\n    [\"17+6=23\", \"67+17=84\"],\n    [\"62-4=58\", \"9+39=48\"],\n    [\"9+60=69\", \"91-28=63\"],\n    [\"3+53=56\", \"4+73=77\"],\n    [\"62-0=62\", \"52+40=92\"],\n    [\"58-29=29\", \"9+12=21\"],\n    [\"85-69=16\", \"11+83=94\"],\n    [\"38+6=44\", \"28-1=27\"],\n    [\"85-72=13\", \"20+60=80\"],\n    [\"91-85=6\", \"25+41=66\"],\n    [\"65+27=92\", \"61+36=97\"],\n    [\"30-2=28\", \"0+14=14\"],\n    [\"96+1=97\", \"72-50=22\"],\n    [\"95-70=25\", \"89-34=55\"],\n    [\"33+17=50\", \"59+26=85\"],\n    [\"31+38=69\", \"9+80=89\"],\n    [\"13+51=64\", \"10+27=37\"],\n    [\"92+0=92\", \"29-25=4\"],\n    [\"93-42=51\", \"87-48=39\"],\n    [\"87-39=48\", \"67-31=36\"],\n    [\"50+15=65\", \"11+71=82\"],\n    [\"69-9=60\", \"32-9=23\"],\n    [\"24+2=26\", \"78-37=41\"],\n    [\"46-28=18\", \"52-32=20\"],\n    [\"96-9=87\", \"15+56=71\"],\n    [\"95-32=63\", \"9+61=70\"],\n    [\"70+8=78\", \"56+7=63\"],\n    [\"68-18=50\", \"40+21=61\"],\n    [\"91-87=4\", \"26+4=30\"],\n    [\"75-25=50\", \"31+52=83\"],\n    [\"44+34=78\", \"98-24=74\"],\n    [\"61+3=64\", \"80-25=55\"],\n    [\"42-39=3\", \"92-34=58\"],\n    [\"67-43=24\", \"39+37=76\"],\n    [\"24+11=35\", \"36-31=5\"],\n    [\"80-29=51\", \"29+68=97\"],\n    [\"23-21=2\", \"31-20=11\"],\n    [\"51+46=97\", \"7+29=36\"],\n    [\"39+58=97\", \"13+76=89\"],\n    [\"52-42=10\", \"72-24=48\"],\n    [\"38+56=94\", \"19-7=12\"],\n    [\"85+3=88\", \"38+60=98\"],\n    [\"43+11=54\", \"99-57=42\"],\n    [\"30+29=59\", \"88-42=46\"],\n    [\"77-8=69\", \"8+63=71\"],\n    [\"65+5=70\", \"44-39=5\"],\n    [\"85-22=63\", \"45+28=73\"],\n    [\"65-45=20\", \"12+69=81\"],\n    [\"65+33=98\", \"88+6=94\"],\n    [\"22+39=61\", \"23+29=52\"],\n    [\"1+66=67\", \"68-41=27\"],\n    [\"1+21=22\", \"17+36=53\"],\n    [\"20+3=23\", \"82+2=84\"],\n    [\"66+23=89\", \"90-85=5\"],\n    [\"21+11=32\", \"56-46=10\"],\n    [\"97-38=59\", \"62-10=52\"],\n    [\"74-42=32\", \"66-9=57\"],\n    [\"89-29=60\", \"36-12=24\"],\n    [\"1+42=43\", \"30-26=4\"],\n    [\"59-15=44\", \"90+6=96\"],\n    [\"26+6=32\", \"38-31=7\"],\n    [\"51-19=32\", \"15-9=6\"],\n    [\"90-75=15\", \"74-66=8\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: true,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each arithmetic expression in the table to its new value.\n# Every expression (e.g. \"86-67=19\") appears exactly once in the document,\n# so Find/Replace on the unique old text is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '86-67=19'; New = '57-16=41' },\n    @{ Old = '52+16=68'; New = '41+28=69' },\n    @{ Old = '28+51=79'; New = '49-37=12' },\n    @{ Old = '70+22=92'; New = '60-24=36' },\n    @{ Old = '19-5=14'; New = '24+70=94' },\n    @{ Old = '88-20=68'; New = '83-30=53' },\n    @{ Old = '13-6=7'; New = '16+21=37' },\n    @{ Old = '62-5=57'; New = '87-43=44' },\n    @{ Old = '36+7=43'; New = '15+16=31' },\n    @{ Old = '86-16=70'; New = '66-49=17' },\n    @{ Old = '51-3=48'; New = '16+82=98' },\n    @{ Old = '74-52=22'; New = '4+38=42' },\n    @{ Old = '41+55=96'; New = '68+24=92' },\n    @{ Old = '54-45=9'; New = '34+29=63' },\n    @{ Old = '12-4=8'; New = '52-16=36' },\n    @{ Old = '45+35=80'; New = '33+66=99' },\n    @{ Old = '8+15=23'; New = '96-65=31' },\n    @{ Old = '72-62=10'; New = '4+52=56' },\n    @{ Old = '78-4=74'; New = '98-91=7' },\n    @{ Old = '93-65=28'; New = '41-39=2' },\n    @{ Old = '14-11=3'; New = '13+68=81' },\n    @{ Old = '27+66=93'; New = '71-54=17' },\n    @{ Old = '28+26=54'; New = '89-14=75' },\n    @{ Old = '61+6=67'; New = '18+67=85' },\n    @{ Old = '20+57=77'; New = '99-59=40' },\n    @{ Old = '18+43=61'; New = '43+53=96' },\n    @{ Old = '40-29=11'; New = '13+20=33' },\n    @{ Old = '23-17=6'; New = '79+4=83' },\n    @{ Old = '64-10=54'; New = '34-34=0' },\n    @{ Old = '78+14=92'; New = '37+28=65' },\n    @{ Old = '67+28=95'; New = '20-9=11' },\n    @{ Old = '32-5=27'; New = '48+17=65' },\n    @{ Old = '2+13=15'; New = '33-2=31' },\n    @{ Old = '86-62=24'; New = '6+16=22' },\n    @{ Old = '86-13=73'; New = '23+1=24' },\n    @{ Old = '76-36=40'; New = '85-10=75' },\n    @{ Old = '93-22=71'; New = '8+32=40' },\n    @{ Old = '17+6=23'; New = '67+17=84' },\n    @{ Old = '62-4=58'; New = '9+39=48' },\n    @{ Old = '9+60=69'; New = '91-28=63' },\n    @{ Old = '3+53=56'; New = '4+73=77' },\n    @{ Old = '62-0=62'; New = '52+40=92' },\n    @{ Old = '58-29=29'; New = '9+12=21' },\n    @{ Old = '85-69=16'; New = '11+83=94' },\n    @{ Old = '38+6=44'; New = '28-1=27' },\n    @{ Old = '85-72=13'; New = '20+60=80' },\n    @{ Old = '91-85=6'; New = '25+41=66' },\n    @{ Old = '65+27=92'; New = '61+36=97' },\n    @{ Old = '30-2=28'; New = '0+14=14' },\n    @{ Old = '96+1=97'; New = '72-50=22' },\n    @{ Old = '95-70=25'; New = '89-34=55' },\n    @{ Old = '33+17=50'; New = '59+26=85' },\n    @{ Old = '31+38=69'; New = '9+80=89' },\n    @{ Old = '13+51=64'; New = '10+27=37' },\n    @{ Old = '92+0=92'; New = '29-25=4' },\n    @{ Old = '93-42=51'; New = '87-48=39' },\n    @{ Old = '87-39=48'; New = '67-31=36' },\n    @{ Old = '50+15=65'; New = '11+71=82' },\n    @{ Old = '69-9=60'; New = '32-9=23' },\n    @{ Old = '24+2=26'; New = '78-37=41' },\n    @{ Old = '46-28=18'; New = '52-32=20' },\n    @{ Old = '96-9=87'; New = '15+56=71' },\n    @{ Old = '95-32=63'; New = '9+61=70' },\n    @{ Old = '70+8=78'; New = '56+7=63' },\n    @{ Old = '68-18=50'; New = '40+21=61' },\n    @{ Old = '91-87=4'; New = '26+4=30' },\n    @{ Old = '75-25=50'; New = '31+52=83' },\n    @{ Old = '44+34=78'; New = '98-24=74' },\n    @{ Old = '61+3=64'; New = '80-25=55' },\n    @{ Old = '42-39=3'; New = '92-34=58' },\n    @{ Old = '67-43=24'; New = '39+37=76' },\n    @{ Old = '24+11=35'; New = '36-31=5' },\n    @{ Old = '80-29=51'; New = '29+68=97' },\n    @{ Old = '23-21=2'; New = '31-20=11' },\n    @{ Old = '51+46=97'; New = '7+29=36' },\n    @{ Old = '39+58=97'; New = '13+76=89' },\n    @{ Old = '52-42=10'; New = '72-24=48' },\n    @{ Old = '38+56=94'; New = '19-7=12' },\n    @{ Old = '85+3=88'; New = '38+60=98' },\n    @{ Old = '43+11=54'; New = '99-57=42' },\n    @{ Old = '30+29=59'; New = '88-42=46' },\n    @{ Old = '77-8=69'; New = '8+63=71' },\n    @{ Old = '65+5=70'; New = '44-39=5' },\n    @{ Old = '85-22=63'; New = '45+28=73' },\n    @{ Old = '65-45=20'; New = '12+69=81' },\n    @{ Old = '65+33=98'; New = '88+6=94' },\n    @{ Old = '22+39=61'; New = '23+29=52' },\n    @{ Old = '1+66=67'; New = '68-41=27' },\n    @{ Old = '1+21=22'; New = '17+36=53' },\n    @{ Old = '20+3=23'; New = '82+2=84' },\n    @{ Old = '66+23=89'; New = '90-85=5' },\n    @{ Old = '21+11=32'; New = '56-46=10' },\n    @{ Old = '97-38=59'; New = '62-10=52' },\n    @{ Old = '74-42=32'; New = '66-9=57' },\n    @{ Old = '89-29=60'; New = '36-12=24' },\n    @{ Old = '1+42=43'; New = '30-26=4' },\n    @{ Old = '59-15=44'; New = '90+6=96' },\n    @{ Old = '26+6=32'; New = '38-31=7' },\n    @{ Old = '51-19=32'; New = '15-9=6' },\n    @{ Old = '90-75=15'; New = '74-66=8' }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $result = $find.Execute(\n        $pair.Old,    # FindText\n        $true,        # MatchCase\n        $true,        # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        0,            # Wrap (wdFindStop)\n        $false,       # Format\n        $pair.New,    # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    )\n\n    if (-not $result) {\n        throw \"Text not found: $($pair.Old)\"\n    }\n}\n"}
